$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates describing the refreshed cryptocurrency market data
$changes = @(
    @{Row=2; D='64.966.06'; E='  -1.53%  '},
    @{Row=3; D='3.135.62'; E='  -7.61%  '},
    @{Row=4; E='  +0.01%  '},
    @{Row=5; D='567.82'; E='  -2.34%  '},
    @{Row=6; D='168.56'; E='  -6.27%  '},
    @{Row=7; E='  -0.02%  '},
    @{Row=8; D='0.600'; E='  -3.46%  '},
    @{Row=9; D='3.134.43'; E='  -7.53%  '},
    @{Row=10; E='  -5.87%  '},
    @{Row=11; E='  -6.01%  '},
    @{Row=12; D='0.389'; E='  -5.71%  '},
    @{Row=13; D='3.673.62'; E='  -7.65%  '},
    @{Row=14; D='0.136'; E='  +1.00%  '},
    @{Row=15; D='26.68'; E='  -7.98%  '},
    @{Row=16; D='64.832.26'; E='  -1.77%  '},
    @{Row=17; D='0.0000161'; E='  -6.09%  '},
    @{Row=18; D='3.128.44'; E='  -8.21%  '},
    @{Row=19; D='5.68'; E='  -3.26%  '},
    @{Row=20; D='12.70'; E='  -7.02%  '},
    @{Row=21; D='354.60'; E='  -3.16%  '},
    @{Row=22; D='7.22'},
    @{Row=23; E='  +0.43%  '},
    @{Row=24; D='68.76'; E='  -5.48%  '},
    @{Row=25; D='0.492'; E='  -6.78%  '},
    @{Row=26; D='3.265.66'; E='  -7.78%  '},
    @{Row=27; D='0.0000114'; E='  -8.05%  '},
    @{Row=28; D='9.63'; E='  -1.12%  '},
    @{Row=29; E='  -2.41%  '},
    @{Row=30; E='  -0.10%  '},
    @{Row=31; D='0.999'; E='  -0.03%  '},
    @{Row=32; E='  -4.16%  '},
    @{Row=33; D='21.72'; E='  -6.14%  '},
    @{Row=34; D='5.24'; E='  -8.50%  '},
    @{Row=35; D='6.55'; E='  -6.24%  '},
    @{Row=36; E='  -5.67%  '},
    @{Row=37; D='158.39'; E='  -1.83%  '},
    @{Row=38; E='  -6.77%  '},
    @{Row=39; D='0.826'; E='  -3.71%  '},
    @{Row=40; D='25.97'; E='  -4.21%  '},
    @{Row=41; D='1.75'; E='  -1.64%  '},
    @{Row=42; D='2.640.82'; E='  -1.12%  '},
    @{Row=43; B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='2.40'; E='  -7.94%  '},
    @{Row=44; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='6.07'; E='  -2.01%  '},
    @{Row=45; D='4.14'; E='  -4.75%  '},
    @{Row=46; D='39.32'; E='  -0.88%  '},
    @{Row=47; D='0.0650'; E='  -4.07%  '},
    @{Row=48; D='23.85'; E='  -2.70%  '},
    @{Row=49; D='317.52'; E='  -5.40%  '},
    @{Row=50; D='0.0270'; E='  -5.24%  '},
    @{Row=51; E='  -2.15%  '}
)

foreach ($change in $changes) {
    $row = $change.Row
    if ($change.ContainsKey('B')) { $ws.Cells.Item($row, 2).Value = $change.B }
    if ($change.ContainsKey('C')) { $ws.Cells.Item($row, 3).Value = $change.C }
    if ($change.ContainsKey('D')) { $ws.Cells.Item($row, 4).Value = $change.D }
    if ($change.ContainsKey('E')) { $ws.Cells.Item($row, 5).Value = $change.E }
}
